# Updated cryptos list on Fri Feb 24 10:10:42 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.
# Numeric-looking Price strings are entered with a leading apostrophe so Excel
# keeps them as text (matching the source data's inlineStr cells) instead of
# coercing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.884.87"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.650.71"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'310.70"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.3896"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "'51.45"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'1.341"
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'0.08446"
$ws.Range("D13").Value = "'23.88"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "'6.998"
$ws.Range("E14").Value = "  -3.92%  "
$ws.Range("D15").Value = "'8.007"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'0.00001315"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "1.651.85"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "'0.06981"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'19.51"
$ws.Range("E20").Value = "  -4.75%  "
$ws.Range("D21").Value = "'6.950"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'13.65"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "23.885.87"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").Value = "'2.447"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("D26").Value = "'2.931"
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").Value = "'22.01"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").Value = "'153.12"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").Value = "'5.404"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").Value = "'137.38"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "'7.747"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'2.483"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").Value = "1.829.70"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").Value = "'0.08136"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "'0.9953"
$ws.Range("E35").Value = "  -5.81%  "
$ws.Range("D36").Value = "'6.658"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "'0.02898"
$ws.Range("E37").Value = "  -6.39%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "'10.69"
$ws.Range("E39").Value = "  -4.18%  "
$ws.Range("D40").Value = "'0.09099"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").Value = "'0.7549"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").Value = "'13.51"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").Value = "'1.415"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").Value = "'16.52"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'0.6922"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("E46").Value = "  -3.47%  "
$ws.Range("D47").Value = "'4.107"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "'0.08268"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'133.36"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("D51").Value = "'1.224"
$ws.Range("E51").Value = "  -2.81%  "
